# Add two more tutoring lesson rows to the hour log (rows 14-15), continuing
# the existing Date / Hours table that ends at row 13.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new rows right below the current last entry, inheriting the
# formatting of the row above (date format in A, centered number in B) -
# the same thing Excel does when you keep typing new rows into the log.
$ws.Rows("14:15").Insert(-4121, 0)  # xlShiftDown, xlFormatFromLeftOrAbove

# New lesson dates (serial 43074 = 2017-12-05, serial 43076 = 2017-12-07).
$ws.Range("A14").Value = 43074
$ws.Range("B14").Value = 2

$ws.Range("A15").Value = 43076
$ws.Range("B15").Value = 2

# Move the active selection to C15, matching the author's post-edit cursor spot.
$ws.Range("C15").Select()

$wb.Application.Calculate()
